$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet "overage" -> "overaged" (tab rename)
# ---------------------------------------------------------------------
$wsOveraged = $wb.Worksheets.Item("overage")
$wsOveraged.Name = "overaged"

# ---------------------------------------------------------------------
# 2) Sheet "access": drop " - Binary" suffix from the disruption label
# ---------------------------------------------------------------------
$wsAccess = $wb.Worksheets.Item("access")
$wsAccess.Range("G4").Value = "% of school-aged children whose education was disrupted due to the school being occupied by displaced persons"

# ---------------------------------------------------------------------
# 3) Sheet "overaged": relabel from "overage" wording to "overaged" wording
#    and reword the two age-related indicator labels
# ---------------------------------------------------------------------
$wsOveraged.Range("A2").Value = "overaged"
$wsOveraged.Range("B2").Value = "Analysis of overaged learners "
$wsOveraged.Range("F2").Value = "Overaged learners"
$wsOveraged.Range("G2").Value = "% of school-aged children attending primary school who are at least 2 years above the intended age for their grade"
$wsOveraged.Range("G3").Value = "% of school-aged children attending secondary school who are at least 2 years above the intended age for their grade"

# ---------------------------------------------------------------------
# 4) Sheet "out_of_school": relabel D2 from "Non Access" to "% of OoS children"
# ---------------------------------------------------------------------
$wsOutOfSchool = $wb.Worksheets.Item("out_of_school")
$wsOutOfSchool.Range("D2").Value = "% of OoS children"

# ---------------------------------------------------------------------
# 5) Selections / scroll positions / active sheet
#    Final active tab must be "out_of_school" (index 2), selection D11.
#    Intermediate selections recorded on "access" and "overaged" sheets.
# ---------------------------------------------------------------------
$wsAccess.Activate()
$wsAccess.Range("G5").Select()

$wsOveraged.Activate()
$wsOveraged.Range("G2:G3").Select()

$wsOutOfSchool.Activate()
$wsOutOfSchool.Range("D11").Select()
